# Update the "FB Account Manager" sheet:
#  - Set column D = old column E values (shift), and set column E = new values,
#    effectively moving the "location" number column D->E by one slot and
#    filling in new figures, per the source diff.
#  - Update the saved selection to E20.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New target values for columns C (unchanged), D, E for rows 1-19.
$data = @{
    1  = @{ D = 41; E = 20 }
    2  = @{ D = 54; E = 27 }
    3  = @{ D = 43; E = 20 }
    4  = @{ D = 80; E = 103 }
    6  = @{ D = 45; E = 20 }
    7  = @{ D = 32; E = 15 }
    8  = @{ D = 54; E = 27 }
    9  = @{ D = 33; E = 15 }
    10 = @{ D = 31; E = 15 }
    12 = @{ D = 31; E = 15 }
    14 = @{ D = 37; E = 15 }
    15 = @{ D = 37; E = 15 }
    16 = @{ D = 30; E = 15 }
    17 = @{ D = 37; E = 15 }
    18 = @{ D = 39; E = 15 }
    19 = @{ D = 54; E = 27 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("E$row").Value = $vals.E
}

# Update the active selection shown in the saved file to E20.
$ws.Range("E20").Select()
